$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"
$ws.Range("C8").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("H22").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("H30").Value = 0
